$wb = $excel.ActiveWorkbook

# Sheet ALC, row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4019.7307
$ws.Range("I40").Value = 2127.4285
$ws.Range("K40").Value = 2127.4285
$ws.Range("M40").Value = -1952.4285

# Sheet ALC, row 69
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 6377.161
$ws.Range("I69").Value = 4237.5
$ws.Range("J69").Value = 6694.148
$ws.Range("K69").Value = 12712.5
$ws.Range("L69").Value = 20082.444
$ws.Range("M69").Value = -11838.5
$ws.Range("N69").Value = -21830.444

# Sheet ALC, row 72
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 6377.161
$ws.Range("I72").Value = 4237.5
$ws.Range("J72").Value = 6694.148
$ws.Range("K72").Value = 38137.5
$ws.Range("L72").Value = 60247.332
$ws.Range("M72").Value = -33769.5
$ws.Range("N72").Value = -68983.33199999999

# Sheet ALC, row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 4562.4443
$ws.Range("I86").Value = 3535.9333
$ws.Range("J86").Value = 5845.5835
$ws.Range("K86").Value = 3535.9333
$ws.Range("L86").Value = 5845.5835
$ws.Range("M86").Value = -2412.9333
$ws.Range("N86").Value = -8091.5835

# Sheet ALC, row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 4562.4443
$ws.Range("I89").Value = 3535.9333
$ws.Range("J89").Value = 5845.5835
$ws.Range("K89").Value = 17679.6665
$ws.Range("L89").Value = 29227.9175
$ws.Range("M89").Value = -12063.6665
$ws.Range("N89").Value = -40459.9175

# Sheet ALC, row 94
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 790
$ws.Range("I94").Value = 790
$ws.Range("K94").Value = 790
$ws.Range("M94").Value = -339

# Sheet ARM, row 80
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 29991.5
$ws.Range("I80").Value = 29991.5
$ws.Range("K80").Value = 29991.5
$ws.Range("M80").Value = -28993.5

# Sheet ARM, row 83
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H83").Value = 29991.5
$ws.Range("I83").Value = 29991.5
$ws.Range("K83").Value = 89974.5
$ws.Range("M83").Value = -84982.5

# Sheet ARM, row 95
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H95").Value = 62499
$ws.Range("J95").Value = 62499
$ws.Range("L95").Value = 62499
$ws.Range("N95").Value = -67991

# Sheet ARM, row 127
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H127").Value = 99755
$ws.Range("J127").Value = 99755
$ws.Range("L127").Value = 99755
$ws.Range("N127").Value = -109675

# Sheet BSM, row 140
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 76666.664
$ws.Range("J140").Value = 76666.664
$ws.Range("L140").Value = 76666.664
$ws.Range("N140").Value = -87026.664

# Sheet CRP, row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4638.1763
$ws.Range("I58").Value = 4980.28
$ws.Range("J58").Value = 3687.889
$ws.Range("K58").Value = 4980.28
$ws.Range("L58").Value = 3687.889
$ws.Range("M58").Value = -4777.28
$ws.Range("N58").Value = -4093.889

# Sheet CRP, row 103
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H103").Value = 9013.333000000001
$ws.Range("J103").Value = 20192
$ws.Range("L103").Value = 20192
$ws.Range("N103").Value = -22536

# Sheet CRP, row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 4638.1763
$ws.Range("I136").Value = 4980.28
$ws.Range("J136").Value = 3687.889
$ws.Range("K136").Value = 14940.84
$ws.Range("L136").Value = 11063.667
$ws.Range("M136").Value = -12390.84
$ws.Range("N136").Value = -16163.667

# Sheet CRP, row 141
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 212727.81
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 212727.81
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 212727.81
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -223087.81

# Sheet CUL, row 37
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 56200
$ws.Range("J37").Value = 56200
$ws.Range("L37").Value = 168600
$ws.Range("N37").Value = -168824

# Sheet CUL, row 69
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 2000
$ws.Range("I69").Value = 2000
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 6000
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -5189
$ws.Range("N69").ClearContents()

# Sheet CUL, row 72
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H72").Value = 2000
$ws.Range("I72").Value = 2000
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 18000
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -13944
$ws.Range("N72").ClearContents()

# Sheet CUL, row 139
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 2243.2144
$ws.Range("I139").Value = 2450.5
$ws.Range("J139").Value = 1725
$ws.Range("K139").Value = 7351.5
$ws.Range("L139").Value = 5175
$ws.Range("M139").Value = -2211.5
$ws.Range("N139").Value = -15455

# Sheet GSM, row 94
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H94").Value = 35000
$ws.Range("J94").Value = 35000
$ws.Range("L94").Value = 35000
$ws.Range("N94").Value = -36352

# Sheet GSM, row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 4762967
$ws.Range("I97").Value = 7937174.5
$ws.Range("K97").Value = 7937174.5
$ws.Range("M97").Value = -7936678.5

# Sheet GSM, row 98
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

# Sheet GSM, row 99
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 8326.556
$ws.Range("I99").Value = 8326.556
$ws.Range("K99").Value = 8326.556
$ws.Range("M99").Value = -6080.556

# Sheet GSM, row 100
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H100").Value = 37490.5
$ws.Range("J100").Value = 37490.5
$ws.Range("L100").Value = 37490.5
$ws.Range("N100").Value = -39654.5

# Sheet GSM, row 127
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H127").Value = 84466.75
$ws.Range("J127").Value = 84466.75
$ws.Range("L127").Value = 84466.75
$ws.Range("N127").Value = -94386.75

# Sheet GSM, row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2517.535
$ws.Range("I132").Value = 2348.6667
$ws.Range("J132").Value = 3386
$ws.Range("K132").Value = 7046.000100000001
$ws.Range("L132").Value = 10158
$ws.Range("M132").Value = -4516.000100000001
$ws.Range("N132").Value = -15218

# Sheet LTW, row 6
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()

# Sheet LTW, row 11
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()

# Sheet LTW, row 17
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 11799.667
$ws.Range("I17").Value = 400
$ws.Range("J17").Value = 17499.5
$ws.Range("K17").Value = 400
$ws.Range("L17").Value = 17499.5
$ws.Range("M17").Value = -230
$ws.Range("N17").Value = -17839.5

# Sheet LTW, row 39
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H39").Value = 8735.200000000001
$ws.Range("J39").Value = 8735.200000000001
$ws.Range("L39").Value = 8735.200000000001
$ws.Range("N39").Value = -9655.200000000001

# Sheet LTW, row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3472851.5
$ws.Range("I61").Value = 3968823.8
$ws.Range("K61").Value = 3968823.8
$ws.Range("M61").Value = -3968621.8

# Sheet LTW, row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 3472851.5
$ws.Range("I113").Value = 3968823.8
$ws.Range("K113").Value = 3968823.8
$ws.Range("M113").Value = -3966653.8

# Sheet LTW, row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4891.926
$ws.Range("I132").Value = 4917.2573
$ws.Range("K132").Value = 14751.7719
$ws.Range("M132").Value = -12221.7719

# Sheet LTW, row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 38207.543
$ws.Range("I136").Value = 74228.53999999999
$ws.Range("J136").Value = 5672.4517
$ws.Range("K136").Value = 222685.62
$ws.Range("L136").Value = 17017.3551
$ws.Range("M136").Value = -220135.62
$ws.Range("N136").Value = -22117.3551

# Sheet WVR, row 99
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H99").Value = 73432
$ws.Range("I99").Value = 73432
$ws.Range("K99").Value = 73432
$ws.Range("M99").Value = -70437

# Sheet WVR, row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1142.6364
$ws.Range("I113").Value = 171.6875
$ws.Range("K113").Value = 515.0625
$ws.Range("M113").Value = 1654.9375

# Sheet WVR, row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 32614668
$ws.Range("I132").Value = 37042972
$ws.Range("J132").Value = 2723596.8
$ws.Range("K132").Value = 111128916
$ws.Range("L132").Value = 8170790.399999999
$ws.Range("M132").Value = -111126386
$ws.Range("N132").Value = -8175850.399999999

# Sheet WVR, row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1288.9354
$ws.Range("I136").Value = 1148.6428
$ws.Range("J136").Value = 2598.3333
$ws.Range("K136").Value = 3445.9284
$ws.Range("L136").Value = 7794.999899999999
$ws.Range("M136").Value = -895.9284000000002
$ws.Range("N136").Value = -12894.9999
